$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    12 = -2
    21 = -2
    22 = -2
    23 = 1
    24 = 2
    25 = 6
    26 = 3
    27 = -2
    30 = 3
    34 = -3
    38 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
